# "Socks in the Dark - 1 a, b, c"
#
# 1) Remove the _GoBack bookmark (left over from the previous edit session)
#    and collapse the two trailing empty paragraphs into one.
# 2/3/4) Answer parts (a), (b) and (c) of question 1 under "Socks in the
#    Dark" with blue-colored commentary runs.

$d = $word.ActiveDocument

$BLUE       = 16711680   # RGB(0x00,0x00,0xFF)
$BLUE_3366  = 16737843   # RGB(0x33,0x66,0xFF)

# ---------------------------------------------------------------------
# 1) Drop the stray _GoBack bookmark and the extra blank paragraph that
#    used to sit beneath it.
# ---------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()
$d.Paragraphs.Item(24).Range.Delete()

# ---------------------------------------------------------------------
# Helper text
# ---------------------------------------------------------------------

# 1a) "In your own words."
$p = $d.Paragraphs.Item(28)
$r = $p.Range
$insertStart = $r.End - 1
$text = "Figure out the likelihood of blindly selecting a matching pair of socks out of 5 black pairs, 3 brown pairs and 2 white pairs."
$r.InsertAfter(" " + $text)
$colored = $d.Range($insertStart + 1, $insertStart + 1 + $text.Length)
$colored.Font.Color = $BLUE

# 1b) "What insight can you offer into the problem..."
$p = $d.Paragraphs.Item(29)
$r = $p.Range
$insertStart = $r.End - 1
$text1 = "The statement does not include any information about the texture or length of the socks.  We do not have to assume they are the same exact brands or styles."
$r.InsertAfter(" " + $text1)
$colored1 = $d.Range($insertStart + 1, $insertStart + 1 + $text1.Length)
$colored1.Font.Color = $BLUE

$p = $d.Paragraphs.Item(29)
$r = $p.Range
$insertStart2 = $r.End - 1
$text2 = "  The problem posed does not take into account neurotic people (like me) who match and roll their socks before storing them."
$r.InsertAfter($text2)
$colored2 = $d.Range($insertStart2, $insertStart2 + $text2.Length)
$colored2.Font.Color = $BLUE

# 1c) "What is the overall goal?"
$p = $d.Paragraphs.Item(30)
$r = $p.Range
$insertStart = $r.End - 1
$r.InsertAfter(" ")

$p = $d.Paragraphs.Item(30)
$r = $p.Range
$textStart = $r.End - 1
$text3 = "To find one matching pair in any color and to find 3 matching pairs of each color"
$r.InsertAfter($text3)
$bmRange = $d.Range($textStart, $textStart + $text3.Length)
$bmRange.Font.Color = $BLUE_3366
$d.Bookmarks.Add("_GoBack", $bmRange)

$p = $d.Paragraphs.Item(30)
$r = $p.Range
$r.InsertAfter(".")
